# fix relation flow bug
# Move the "Relation_In_Out"/"Relation_Out_Out" values recorded for the
# Electrolyzer row into the correct "Relation_In_In" column, and correct the
# Methanol_Reactor row's Relation_Out_Out value.

$wb = $excel.ActiveWorkbook
$wsUnits = $wb.Worksheets.Item("Units")
$wsConnections = $wb.Worksheets.Item("Connections")

# Row 3 (Electrolyzer): the value that was wrongly entered under
# Relation_In_Out (P3) actually belongs under Relation_In_In (O3); the
# spurious Relation_Out_Out entry (Q3) is removed entirely.
$wsUnits.Range("O3").Value = 2
$wsUnits.Range("P3").ClearContents()
$wsUnits.Range("Q3").ClearContents()

# Row 6 (Methanol_Reactor): correct the Relation_Out_Out value.
$wsUnits.Range("Q6").Value = 4

# Update sheet selections/active view to match where the fix was made.
$wsConnections.Activate()
$wsConnections.Range("F11").Select()

$wsUnits.Activate()
$wsUnits.Range("O4").Select()
